# issue #5: stock data from json to db
# Adds "category", "source_file" and "index" columns (metadata describing
# which scraped json file each row came from) to the 股票 (stock) worksheet,
# and fixes a bad face_value row (E13: 1 -> 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Last data row on this sheet.
$lastRow = 15

# Insert a new column at I (pushes old date/legislator_name/legislator_id
# from I/J/K to J/K/L, copying their formatting).
$ws.Columns("I:I").Insert()

# New column I: "category" - the source file's category (the "normal" vs.
# "spouse"/... bucket the report json came from), constant for this sheet.
$ws.Cells.Item(1, 9).Value = "category"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
}

# Insert 2 more new trailing columns, M (source_file) and N (index), right
# after L (legislator_id) so they inherit L's header/data formatting.
$ws.Columns("M:N").Insert()

$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = "tmped981"
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value2
}

# Fix bad face_value for row 13 (was 1, should be 10).
$ws.Cells.Item(13, 5).Value = 10
